$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.523.96"
$ws.Range("E2").Value = "  +1.04%  "

# Row 3
$ws.Range("D3").Value = "1.727.44"

# Row 4
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'244.33"
$ws.Range("E5").Value = "  +2.00%  "

# Row 6
$ws.Range("D6").Value = "'0.9997"

# Row 7
$ws.Range("D7").Value = "'0.4808"
$ws.Range("E7").Value = "  +1.92%  "

# Row 8
$ws.Range("D8").Value = "'0.2668"
$ws.Range("E8").Value = "  +1.85%  "

# Row 9
$ws.Range("D9").Value = "'0.06185"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("D10").Value = "1.732.13"
$ws.Range("E10").Value = "  +0.86%  "

# Row 11
$ws.Range("D11").Value = "'0.07184"
$ws.Range("E11").Value = "  +1.56%  "

# Row 12
$ws.Range("E12").Value = "  +1.19%  "

# Row 13
$ws.Range("D13").Value = "'0.6108"
$ws.Range("E13").Value = "  +2.58%  "

# Row 14
$ws.Range("D14").Value = "'4.530"
$ws.Range("E14").Value = "  +2.43%  "

# Row 15
$ws.Range("D15").Value = "'77.14"
$ws.Range("E15").Value = "  +1.30%  "

# Row 16
$ws.Range("E16").Value = "  -0.08%  "

# Row 17
$ws.Range("D17").Value = "26.537.10"
$ws.Range("E17").Value = "  +1.05%  "

# Row 18
$ws.Range("D18").Value = "'0.9998"
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("D19").Value = "'0.000006932"
$ws.Range("E19").Value = "  +1.98%  "

# Row 20
$ws.Range("D20").Value = "'11.53"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("D21").Value = "1.956.98"
$ws.Range("E21").Value = "  +1.00%  "

# Row 22
$ws.Range("D22").Value = "'4.528"
$ws.Range("E22").Value = "  -0.19%  "

# Row 23
$ws.Range("D23").Value = "'8.796"
$ws.Range("E23").Value = "  +0.82%  "

# Row 24
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").Value = "'137.09"
$ws.Range("E25").Value = "  +1.67%  "

# Row 26
$ws.Range("D26").Value = "'15.33"
$ws.Range("E26").Value = "  +1.14%  "

# Row 27
$ws.Range("D27").Value = "'1.779"

# Row 28
$ws.Range("D28").Value = "'1.398"
$ws.Range("E28").Value = "  -0.27%  "

# Row 29
$ws.Range("D29").Value = "'107.26"
$ws.Range("E29").Value = "  +0.31%  "

# Row 30
$ws.Range("D30").Value = "'3.966"
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("D31").Value = "'0.08024"
$ws.Range("E31").Value = "  +3.50%  "

# Row 32
$ws.Range("D32").Value = "'3.692"
$ws.Range("E32").Value = "  +0.36%  "

# Row 33
$ws.Range("D33").Value = "'0.04521"
$ws.Range("E33").Value = "  +1.16%  "

# Row 34
$ws.Range("E34").Value = "  -0.15%  "

# Row 35
$ws.Range("D35").Value = "'0.9976"
$ws.Range("E35").Value = "  +2.33%  "

# Row 36
$ws.Range("D36").Value = "'0.6264"
$ws.Range("E36").Value = "  +1.49%  "

# Row 37
$ws.Range("D37").Value = "'0.9122"
$ws.Range("E37").Value = "  -1.27%  "

# Row 38
$ws.Range("D38").Value = "'2.073"
$ws.Range("E38").Value = "  +8.00%  "

# Row 39
$ws.Range("D39").Value = "'2.371"
$ws.Range("E39").Value = "  -2.61%  "

# Row 40
$ws.Range("E40").Value = "  -0.01%  "

# Row 41
$ws.Range("E41").Value = "  -9.30%  "

# Row 42
$ws.Range("D42").Value = "'0.01503"
$ws.Range("E42").Value = "  +1.45%  "

# Row 43
$ws.Range("D43").Value = "'5.628"
$ws.Range("E43").Value = "  +0.57%  "

# Row 44
$ws.Range("E44").Value = "  +1.20%  "

# Row 45
$ws.Range("D45").Value = "'6.976"
$ws.Range("E45").Value = "  +10.99%  "

# Row 46
$ws.Range("D46").Value = "'0.1181"
$ws.Range("E46").Value = "  +0.39%  "

# Row 47
$ws.Range("D47").Value = "'0.05364"
$ws.Range("E47").Value = "  +1.81%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'30.48"
$ws.Range("E48").Value = "  +0.38%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.771"
$ws.Range("E49").Value = "  +0.15%  "

# Row 50
$ws.Range("D50").Value = "'1.254"
$ws.Range("E50").Value = "  +3.10%  "

# Row 51
$ws.Range("D51").Value = "'51.33"
$ws.Range("E51").Value = "  +1.39%  "
